$d = $word.ActiveDocument
$d.Content.Find.Execute("SOEN 8641", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SOEN 6841", 2)
